# Scheduled-runner refresh of market/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-class Leve sheets. Static cached
# values only (no formulas in these columns) -- just overwrite with the
# freshly scraped numbers.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4114.2856
$ws.Range("I10").Value = 3950
$ws.Range("J10").Value = 4333.3335
$ws.Range("K10").Value = 3950
$ws.Range("L10").Value = 4333.3335
$ws.Range("M10").Value = -3657
$ws.Range("N10").Value = -4919.3335
$ws.Range("H13").Value = 14950
$ws.Range("J13").Value = 19900
$ws.Range("L13").Value = 19900
$ws.Range("N13").Value = -20238
$ws.Range("H43").Value = 490.07144
$ws.Range("I43").Value = 275.4
$ws.Range("J43").Value = 609.3333
$ws.Range("K43").Value = 275.4
$ws.Range("L43").Value = 609.3333
$ws.Range("M43").Value = -206.4
$ws.Range("N43").Value = -747.3333
$ws.Range("H64").Value = 3887.15
$ws.Range("I64").Value = 3595.5557
$ws.Range("J64").Value = 4125.727
$ws.Range("K64").Value = 3595.5557
$ws.Range("L64").Value = 4125.727
$ws.Range("M64").Value = -3347.5557
$ws.Range("N64").Value = -4621.727
$ws.Range("H67").Value = 3887.15
$ws.Range("I67").Value = 3595.5557
$ws.Range("J67").Value = 4125.727
$ws.Range("K67").Value = 3595.5557
$ws.Range("L67").Value = 4125.727
$ws.Range("M67").Value = -2737.5557
$ws.Range("N67").Value = -5841.727
$ws.Range("H112").Value = 1702.6364
$ws.Range("I112").Value = 345
$ws.Range("J112").Value = 2004.3334
$ws.Range("K112").Value = 1035
$ws.Range("L112").Value = 6013.0002
$ws.Range("M112").Value = 73
$ws.Range("N112").Value = -8229.0002
$ws.Range("H116").Value = 2199.7896
$ws.Range("J116").Value = 2345.0908
$ws.Range("L116").Value = 2345.0908
$ws.Range("N116").Value = -9229.0908
$ws.Range("H137").Value = 2487.025
$ws.Range("I137").Value = 2225.1936
$ws.Range("J137").Value = 3388.889
$ws.Range("K137").Value = 6675.5808
$ws.Range("L137").Value = 10166.667
$ws.Range("M137").Value = -4125.5808
$ws.Range("N137").Value = -15266.667
$ws.Range("H138").Value = 4091.691
$ws.Range("I138").Value = 535.37933
$ws.Range("J138").Value = 8058.346
$ws.Range("K138").Value = 1606.13799
$ws.Range("L138").Value = 24175.038
$ws.Range("M138").Value = 3533.86201
$ws.Range("N138").Value = -34455.038
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28244.72
$ws.Range("I32").Value = 32359.117
$ws.Range("K32").Value = 32359.117
$ws.Range("M32").Value = -32072.117
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26802
$ws.Range("H132").Value = 2467.025
$ws.Range("I132").Value = 1645.1818
$ws.Range("J132").Value = 3471.5
$ws.Range("K132").Value = 4935.5454
$ws.Range("L132").Value = 10414.5
$ws.Range("M132").Value = -2405.5454
$ws.Range("N132").Value = -15474.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1060
$ws.Range("I20").Value = 1017.6923
$ws.Range("J20").Value = 1197.5
$ws.Range("K20").Value = 1017.6923
$ws.Range("L20").Value = 1197.5
$ws.Range("M20").Value = -770.6923
$ws.Range("N20").Value = -1691.5
$ws.Range("H94").Value = 2300
$ws.Range("I94").Value = 2300
$ws.Range("K94").Value = 2300
$ws.Range("M94").Value = -1849
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3033632.8
$ws.Range("I58").Value = 8267517
$ws.Range("J58").Value = 3489
$ws.Range("K58").Value = 8267517
$ws.Range("L58").Value = 3489
$ws.Range("M58").Value = -8267314
$ws.Range("N58").Value = -3895
$ws.Range("H134").Value = 3572
$ws.Range("I134").Value = 2150.5
$ws.Range("K134").Value = 6451.5
$ws.Range("M134").Value = -3916.5
$ws.Range("H136").Value = 3033632.8
$ws.Range("I136").Value = 8267517
$ws.Range("J136").Value = 3489
$ws.Range("K136").Value = 24802551
$ws.Range("L136").Value = 10467
$ws.Range("M136").Value = -24800001
$ws.Range("N136").Value = -15567
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 820
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -188
$ws.Range("H5").Value = 11912438
$ws.Range("I5").Value = 614.375
$ws.Range("J5").Value = 27794870
$ws.Range("K5").Value = 1843.125
$ws.Range("L5").Value = 83384610
$ws.Range("M5").Value = -1731.125
$ws.Range("N5").Value = -83384834
$ws.Range("H9").Value = 34862.07
$ws.Range("I9").Value = 1500
$ws.Range("J9").Value = 37333.332
$ws.Range("K9").Value = 4500
$ws.Range("L9").Value = 111999.996
$ws.Range("M9").Value = -4276
$ws.Range("N9").Value = -112447.996
$ws.Range("H20").Value = 1320
$ws.Range("H63").Value = 3642.5
$ws.Range("I63").Value = 1749
$ws.Range("J63").Value = 3958.0833
$ws.Range("K63").Value = 5247
$ws.Range("L63").Value = 11874.2499
$ws.Range("M63").Value = -4498
$ws.Range("N63").Value = -13372.2499
$ws.Range("H66").Value = 3642.5
$ws.Range("I66").Value = 1749
$ws.Range("J66").Value = 3958.0833
$ws.Range("K66").Value = 15741
$ws.Range("L66").Value = 35622.7497
$ws.Range("M66").Value = -11997
$ws.Range("N66").Value = -43110.7497
$ws.Range("H122").Value = 788.03845
$ws.Range("I122").Value = 324.16666
$ws.Range("J122").Value = 927.2
$ws.Range("K122").Value = 2917.49994
$ws.Range("L122").Value = 8344.800000000001
$ws.Range("M122").Value = -467.4999399999997
$ws.Range("N122").Value = -13244.8
$ws.Range("H135").Value = 11912438
$ws.Range("I135").Value = 614.375
$ws.Range("J135").Value = 27794870
$ws.Range("K135").Value = 5529.375
$ws.Range("L135").Value = 250153830
$ws.Range("M135").Value = -2994.375
$ws.Range("N135").Value = -250158900
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 34.4
$ws.Range("I2").Value = 37
$ws.Range("J2").Value = 31.8
$ws.Range("K2").Value = 37
$ws.Range("L2").Value = 31.8
$ws.Range("M2").Value = 76
$ws.Range("N2").Value = -257.8
$ws.Range("H31").Value = 1065.5
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 1065.5
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H102").Value = 3137.1353
$ws.Range("I102").Value = 3148.4348
$ws.Range("J102").Value = 3118.5715
$ws.Range("K102").Value = 3148.4348
$ws.Range("L102").Value = 3118.5715
$ws.Range("M102").Value = -1526.4348
$ws.Range("N102").Value = -6362.5715
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2152.1538
$ws.Range("I46").Value = 2183
$ws.Range("J46").Value = 2116.1667
$ws.Range("K46").Value = 2183
$ws.Range("L46").Value = 2116.1667
$ws.Range("M46").Value = -1995
$ws.Range("N46").Value = -2492.1667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 2996.6667
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5280
$ws.Range("H62").Value = 3719.3684
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 3820.6155
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 3820.6155
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -5068.6155
$ws.Range("H65").Value = 3719.3684
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 3820.6155
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 19103.0775
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -25343.0775
$ws.Range("H126").Value = 1701.5
$ws.Range("I126").Value = 1651.8
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 4955.4
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -2485.4
$ws.Range("N126").Value = -10790
$ws.Range("H132").Value = 1791.75
$ws.Range("I132").Value = 698.2917
$ws.Range("J132").Value = 3103.9
$ws.Range("K132").Value = 2094.8751
$ws.Range("L132").Value = 9311.700000000001
$ws.Range("M132").Value = 435.1248999999998
$ws.Range("N132").Value = -14371.7
